$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "Última actualización: 14:17:27"
$ws1.Range("A3").Value = "Total filas: 307"
$ws1.Range("A61").Value = "05:47:32"
$ws1.Range("C61").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D61").Value = 105
$ws1.Range("A62").Value = "06:02:16"
$ws1.Range("C62").Value = "11_ETCHEVERRY"
$ws1.Range("D62").Value = 90
$ws1.Range("A91").Value = "08:47:19"
$ws1.Range("C91").Value = "10_OLMOS"
$ws1.Range("D91").Value = 7
$ws1.Range("A92").Value = "07:14:27"
$ws1.Range("C92").Value = "17_ROMERO"
$ws1.Range("D92").Value = 100
$ws1.Range("A111").Value = "07:57:27"
$ws1.Range("C111").Value = "11_ETCHEVERRY"
$ws1.Range("D111").Value = 86
$ws1.Range("A112").Value = "08:47:19"
$ws1.Range("C112").Value = "16_SANTA ANA"
$ws1.Range("D112").Value = 36
$ws1.Range("A113").Value = "07:44:08"
$ws1.Range("C113").Value = "17_ROMERO"
$ws1.Range("D113").Value = 99
$ws1.Range("A118").Value = "08:33:47"
$ws1.Range("C118").Value = "16_SANTA ANA"
$ws1.Range("D118").Value = 61
$ws1.Range("A119").Value = "08:54:42"
$ws1.Range("C119").Value = "23_HERNANDEZ"
$ws1.Range("D119").Value = 40
$ws1.Range("C120").Value = "16_SANTA ANA"
$ws1.Range("C121").Value = "23_HERNANDEZ"
$ws1.Range("A140").Value = "10:11:11"
$ws1.Range("C140").Value = "16_SANTA ANA"
$ws1.Range("D140").Value = 12
$ws1.Range("A141").Value = "09:25:30"
$ws1.Range("C141").Value = "11_ETCHEVERRY"
$ws1.Range("D141").Value = 58
$ws1.Range("A178").Value = "11:15:53"
$ws1.Range("C178").Value = "23_HERNANDEZ"
$ws1.Range("D178").Value = 20
$ws1.Range("A179").Value = "10:50:37"
$ws1.Range("C179").Value = "11_ETCHEVERRY"
$ws1.Range("D179").Value = 45
$ws1.Range("A188").Value = "11:15:53"
$ws1.Range("C188").Value = "225_GOMEZ"
$ws1.Range("D188").Value = 43
$ws1.Range("A189").Value = "11:58:46"
$ws1.Range("C189").Value = "17_ROMERO"
$ws1.Range("D189").Value = 0
$ws1.Range("C197").Value = "14_ABASTO"
$ws1.Range("A198").Value = "10:11:11"
$ws1.Range("C198").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D198").Value = 116
$ws1.Range("A199").Value = "10:50:37"
$ws1.Range("C199").Value = "10_OLMOS"
$ws1.Range("D199").Value = 77
$ws1.Range("C217").Value = "23_HERNANDEZ"
$ws1.Range("C218").Value = "27_EL RETIRO"
$ws1.Range("C264").Value = "16_P MOR-167 Y 521"
$ws1.Range("C265").Value = "225_GOMEZ"
$ws1.Range("A278").Value = "14:17:27"
$ws1.Range("B278").Value = "14:18"
$ws1.Range("C278").Value = "16_SANTA ANA"
$ws1.Range("D278").Value = 1
$ws1.Range("B279").Value = "14:19"
$ws1.Range("C279").Value = "215C_EL PATO"
$ws1.Range("D279").Value = 28
$ws1.Range("A280").Value = "13:51:56"
$ws1.Range("C280").Value = "26_HERNANDEZ"
$ws1.Range("D280").Value = 29
$ws1.Range("B281").Value = "14:20"
$ws1.Range("C281").Value = "215C_EL PATO"
$ws1.Range("D281").Value = 116
$ws1.Range("A282").Value = "12:24:14"
$ws1.Range("B282").Value = "14:21"
$ws1.Range("C282").Value = "26_HERNANDEZ"
$ws1.Range("D282").Value = 117
$ws1.Range("A283").Value = "14:17:27"
$ws1.Range("B283").Value = "14:28"
$ws1.Range("C283").Value = "15_ABASTO"
$ws1.Range("D283").Value = 11
$ws1.Range("A284").Value = "14:17:27"
$ws1.Range("B284").Value = "14:29"
$ws1.Range("C284").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D284").Value = 12
$ws1.Range("A285").Value = "14:17:27"
$ws1.Range("B285").Value = "14:30"
$ws1.Range("C285").Value = "16_SANTA ANA"
$ws1.Range("D285").Value = 13
$ws1.Range("A286").Value = "14:17:27"
$ws1.Range("B286").Value = "14:34"
$ws1.Range("C286").Value = "23_HERNANDEZ"
$ws1.Range("D286").Value = 17
$ws1.Range("A287").Value = "13:51:56"
$ws1.Range("B287").Value = "14:44"
$ws1.Range("C287").Value = "10_OLMOS"
$ws1.Range("D287").Value = 53
$ws1.Range("A288").Value = "13:51:56"
$ws1.Range("B288").Value = "14:44"
$ws1.Range("C288").Value = "14_ABASTO"
$ws1.Range("D288").Value = 53
$ws1.Range("A289").Value = "12:57:33"
$ws1.Range("B289").Value = "14:45"
$ws1.Range("C289").Value = "14_ABASTO"
$ws1.Range("D289").Value = 108
$ws1.Range("A290").Value = "12:57:33"
$ws1.Range("B290").Value = "14:56"
$ws1.Range("C290").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D290").Value = 119
$ws1.Range("B291").Value = "14:57"
$ws1.Range("C291").Value = "215B_EL PATO"
$ws1.Range("D291").Value = 66
$ws1.Range("B292").Value = "14:58"
$ws1.Range("C292").Value = "215B_EL PATO"
$ws1.Range("D292").Value = 94
$ws1.Range("A293").Value = "14:17:27"
$ws1.Range("B293").Value = "15:00"
$ws1.Range("C293").Value = "10_OLMOS"
$ws1.Range("D293").Value = 43
$ws1.Range("B294").Value = "15:00"
$ws1.Range("C294").Value = "81_EL PELIGRO"
$ws1.Range("D294").Value = 96
$ws1.Range("B295").Value = "15:04"
$ws1.Range("C295").Value = "10_OLMOS"
$ws1.Range("D295").Value = 73
$ws1.Range("A296").Value = "13:24:27"
$ws1.Range("B296").Value = "15:05"
$ws1.Range("C296").Value = "10_OLMOS"
$ws1.Range("D296").Value = 101
$ws1.Range("A297").Value = "14:17:27"
$ws1.Range("B297").Value = "15:10"
$ws1.Range("C297").Value = "17_ROMERO"
$ws1.Range("D297").Value = 53
$ws1.Range("B298").Value = "15:13"
$ws1.Range("C298").Value = "11_ETCHEVERRY"
$ws1.Range("D298").Value = 82
$ws1.Range("A299").Value = "13:24:27"
$ws1.Range("B299").Value = "15:20"
$ws1.Range("C299").Value = "15_ABASTO"
$ws1.Range("D299").Value = 116
$ws1.Range("A300").Value = "13:51:56"
$ws1.Range("B300").Value = "15:21"
$ws1.Range("C300").Value = "26_HERNANDEZ"
$ws1.Range("D300").Value = 90
$ws1.Range("E300").Value = "LP1912"
$ws1.Range("A301").Value = "13:24:27"
$ws1.Range("B301").Value = "15:22"
$ws1.Range("C301").Value = "26_HERNANDEZ"
$ws1.Range("D301").Value = 118
$ws1.Range("E301").Value = "LP1912"
$ws1.Range("A302").Value = "13:51:56"
$ws1.Range("B302").Value = "15:31"
$ws1.Range("C302").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D302").Value = 100
$ws1.Range("E302").Value = "LP1912"
$ws1.Range("A303").Value = "14:17:27"
$ws1.Range("B303").Value = "15:32"
$ws1.Range("C303").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D303").Value = 75
$ws1.Range("E303").Value = "LP1912"
$ws1.Range("A304").Value = "13:51:56"
$ws1.Range("B304").Value = "15:34"
$ws1.Range("C304").Value = "23_HERNANDEZ"
$ws1.Range("D304").Value = 103
$ws1.Range("E304").Value = "LP1912"
$ws1.Range("A305").Value = "13:51:56"
$ws1.Range("B305").Value = "15:38"
$ws1.Range("C305").Value = "215A_EL PATO"
$ws1.Range("D305").Value = 107
$ws1.Range("E305").Value = "LP1912"
$ws1.Range("A306").Value = "14:17:27"
$ws1.Range("B306").Value = "15:38"
$ws1.Range("C306").Value = "23_HERNANDEZ"
$ws1.Range("D306").Value = 81
$ws1.Range("E306").Value = "LP1912"
$ws1.Range("A307").Value = "13:51:56"
$ws1.Range("B307").Value = "15:46"
$ws1.Range("C307").Value = "16_P MOR-167 Y 521"
$ws1.Range("D307").Value = 115
$ws1.Range("E307").Value = "LP1912"
$ws1.Range("A308").Value = "13:51:56"
$ws1.Range("B308").Value = "15:49"
$ws1.Range("C308").Value = "14_ABASTO"
$ws1.Range("D308").Value = 118
$ws1.Range("E308").Value = "LP1912"
$ws1.Range("A309").Value = "14:17:27"
$ws1.Range("B309").Value = "15:56"
$ws1.Range("C309").Value = "27_EL RETIRO"
$ws1.Range("D309").Value = 99
$ws1.Range("E309").Value = "LP1912"
$ws1.Range("A310").Value = "14:17:27"
$ws1.Range("B310").Value = "15:56"
$ws1.Range("C310").Value = "17_ROMERO"
$ws1.Range("D310").Value = 99
$ws1.Range("E310").Value = "LP1912"
$ws1.Range("A311").Value = "14:17:27"
$ws1.Range("B311").Value = "15:57"
$ws1.Range("C311").Value = "11_ETCHEVERRY"
$ws1.Range("D311").Value = 100
$ws1.Range("E311").Value = "LP1912"
$ws1.Range("A312").Value = "14:17:27"
$ws1.Range("B312").Value = "16:15"
$ws1.Range("C312").Value = "225_C ROCA-H SUR"
$ws1.Range("D312").Value = 118
$ws1.Range("E312").Value = "LP1912"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Última actualización: 14:17:27"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "Última actualización: 14:17:27"
$ws3.Range("A3").Value = "Total filas: 42"
$ws3.Range("A46").Value = "14:17:27"
$ws3.Range("B46").Value = "15:34"
$ws3.Range("C46").Value = "215A_LA PLATA"
$ws3.Range("D46").Value = 77
$ws3.Range("E46").Value = "L6173"
$ws3.Range("A47").Value = "14:17:27"
$ws3.Range("B47").Value = "16:14"
$ws3.Range("C47").Value = "215C_LA PLATA"
$ws3.Range("D47").Value = 117
$ws3.Range("E47").Value = "L6203"
